$d = $word.ActiveDocument

# The document has one section whose header/footer "first page" and
# "primary" stories each carry one inline picture:
#   - Header(2)  -> BTec logo jpg   : name "image2.jpg" -> "image1.jpg"
#   - Footer(1)  -> Pearson logo png: name "image1.png" -> "image2.png"
#   - Footer(2)  -> Pearson logo png: name "image1.png" -> "image2.png"
$section = $d.Sections(1)

$btecLogo = $section.Headers(2).Range.InlineShapes(1)
$btecLogo.Name = "image1.jpg"

$pearsonLogoA = $section.Footers(1).Range.InlineShapes(1)
$pearsonLogoA.Name = "image2.png"

$pearsonLogoB = $section.Footers(2).Range.InlineShapes(1)
$pearsonLogoB.Name = "image2.png"
